$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add totals formulas across the "TOTAL" row (row 8), columns D through V
$ws.Range("D8").Formula = "=SUM(D5:D7)"
$ws.Range("E8").Formula = "=(SUM(E5:E7))"
$ws.Range("F8:V8").FormulaR1C1 = "=(SUM(R[-3]C:R[-1]C))"

# Move the active selection to E8, matching the saved selection state
$ws.Range("E8").Select()
